# C5-PowerPoint.pptx — apply the "Table Design" gallery style change made to
# the sources-of-finance table on slide 6 (Table_0 -> built-in
# "Medium Style 2 - Accent 1" table style).
#
# NB: the original commit also toggled File > Options > Save > "Embed fonts
# in the file" (which rewrites ppt/presentation.xml's embedTrueTypeFonts
# attribute, adds a p:embeddedFontLst pointing at an embedded Limelight.ttf,
# and re-applied the built-in "Office Theme" design, which causes PowerPoint
# to shuffle the Integral/Office Theme content between theme1.xml/theme2.xml
# on save). Those are save-dialog/gallery-only operations with no surface in
# the Presentation COM object model (Presentation.Fonts/Font.Embedded are
# read-only, there is no "embed font" method, and ApplyTheme needs a real
# .thmx file on disk) so they cannot be scripted here.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(6)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{99DF9202-A60F-44C5-A9E4-B42114836971}")
    }
}
